# Update cryptocurrency price/volume data as scraped on 2024-04-17
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell "D2" "63.536.06"
Set-TextCell "E2" "  +0.26%  "
Set-TextCell "D3" "3.083.90"
Set-TextCell "E3" "  -0.38%  "
Set-TextCell "D4" "1.00"
Set-TextCell "E4" "  +0.07%  "
Set-TextCell "D5" "545.23"
Set-TextCell "E5" "  -0.84%  "
Set-TextCell "D6" "139.70"
Set-TextCell "E6" "  +1.67%  "
Set-TextCell "D8" "3.077.26"
Set-TextCell "E8" "  -0.36%  "
Set-TextCell "D9" "0.499"
Set-TextCell "E9" "  +0.33%  "
Set-TextCell "D10" "0.157"
Set-TextCell "E10" "  +0.61%  "
Set-TextCell "E11" "  +2.57%  "
Set-TextCell "D12" "0.458"
Set-TextCell "E12" "  -2.69%  "
Set-TextCell "B13" "ShibaInu"
Set-TextCell "C13" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextCell "D13" "0.0000225"
Set-TextCell "E13" "  +3.71%  "
Set-TextCell "B14" "Avalanche"
Set-TextCell "C14" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell "D14" "35.12"
Set-TextCell "E14" "  -1.36%  "
Set-TextCell "D15" "3.583.54"
Set-TextCell "E15" "  -0.14%  "
Set-TextCell "D16" "63.555.50"
Set-TextCell "E16" "  +0.41%  "
Set-TextCell "E17" "  +1.07%  "
Set-TextCell "D18" "3.080.76"
Set-TextCell "E18" "  -0.48%  "
Set-TextCell "E19" "  -1.28%  "
Set-TextCell "D20" "476.12"
Set-TextCell "E20" "  -2.75%  "
Set-TextCell "D21" "13.51"
Set-TextCell "E21" "  -1.22%  "
Set-TextCell "E22" "  -2.65%  "
Set-TextCell "E23" "  -2.37%  "
Set-TextCell "D24" "78.75"
Set-TextCell "E24" "  -0.63%  "
Set-TextCell "D25" "12.27"
Set-TextCell "E25" "  -1.13%  "
Set-TextCell "E26" "  +0.29%  "
Set-TextCell "E27" "  -1.16%  "
Set-TextCell "D28" "7.98"
Set-TextCell "E28" "  -5.85%  "
Set-TextCell "D29" "0.999"
Set-TextCell "E29" "  -0.05%  "
Set-TextCell "D30" "26.29"
Set-TextCell "E30" "  -1.34%  "
Set-TextCell "D31" "1.90"
Set-TextCell "E31" "  -3.74%  "
Set-TextCell "E32" "  +1.80%  "
Set-TextCell "D33" "58.24"
Set-TextCell "E33" "  +0.33%  "
Set-TextCell "E34" "  -7.44%  "
Set-TextCell "D35" "5.49"
Set-TextCell "E35" "  +6.55%  "
Set-TextCell "D36" "492.79"
Set-TextCell "E36" "  -4.35%  "
Set-TextCell "E37" "  -0.06%  "
Set-TextCell "D38" "3.264.00"
Set-TextCell "E38" "  +3.33%  "
Set-TextCell "E39" "  +0.51%  "
Set-TextCell "E40" "  -0.47%  "
Set-TextCell "E41" "  -1.68%  "
Set-TextCell "E42" "  -0.16%  "
Set-TextCell "D43" "2.61"
Set-TextCell "E43" "  -2.08%  "
Set-TextCell "E44" "  -1.97%  "
Set-TextCell "E45" "  +0.05%  "
Set-TextCell "D46" "25.63"
Set-TextCell "E46" "  +1.34%  "
Set-TextCell "D47" "123.48"
Set-TextCell "E47" "  +2.10%  "
Set-TextCell "D48" "2.03"
Set-TextCell "E48" "  -1.67%  "
Set-TextCell "E49" "  +5.49%  "
Set-TextCell "D50" "0.109"
Set-TextCell "E50" "  +0.70%  "
Set-TextCell "E51" "  -0.57%  "
